# Generate Report for Handoff
# Adds two new localization entries:
#   - 6d2a7d87-f970-42ce-bf51-6d3c850207dd.md  (inserted as row 3, before cbcdd771)
#   - edbab8bf-23e4-47e8-831b-14c9aee13e0e.md  (appended as row 5 / last)
# cbcdd771-a78d-4d08-b66d-488e5b202a24.md shifts down from row 3 to row 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 3 now belongs to the new 6d2a7d87 file (was cbcdd771).
$ws.Range("A3").Value = "6d2a7d87-f970-42ce-bf51-6d3c850207dd.md"
$ws.Range("B3").Value = "e2e\6d2a7d87-f970-42ce-bf51-6d3c850207dd.md"
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-21 00:47:53"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 4: cbcdd771 (moved down from row 3).
$ws.Range("A4").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$ws.Range("B4").Value = "e2e\cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$ws.Range("C4").Value = ".md"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-21 00:46:51"
$ws.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 5: new edbab8bf file.
$ws.Range("A5").Value = "edbab8bf-23e4-47e8-831b-14c9aee13e0e.md"
$ws.Range("B5").Value = "e2e\edbab8bf-23e4-47e8-831b-14c9aee13e0e.md"
$ws.Range("C5").Value = ".md"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-21 00:47:53"
$ws.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks in column B (Range.Hyperlinks.Delete() clears the whole
# sheet's collection, so re-add every row in ref order to keep rId2..rId5
# matching the intended sequence).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a454c3bc198b02c2e49d4016dd6550af4a97165/e2e/98ad48ff-e5a0-4540-a055-8ec88ce579da.md", "", "", "e2e\98ad48ff-e5a0-4540-a055-8ec88ce579da.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197eded308f8ed22b21c9850f80743b3961c023e/e2e/6d2a7d87-f970-42ce-bf51-6d3c850207dd.md", "", "", "e2e\6d2a7d87-f970-42ce-bf51-6d3c850207dd.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb043c5240aa15b0c64d9bbc34b532dce3e332d3/e2e/cbcdd771-a78d-4d08-b66d-488e5b202a24.md", "", "", "e2e\cbcdd771-a78d-4d08-b66d-488e5b202a24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/977c565182eb70f5592892745c773df547f5cf39/e2e/edbab8bf-23e4-47e8-831b-14c9aee13e0e.md", "", "", "e2e\edbab8bf-23e4-47e8-831b-14c9aee13e0e.md") | Out-Null

# Grow the "Overview" table (and its autofilter) to cover the new rows.
$loOverview = $ws.ListObjects.Item(1)
$loOverview.Resize($ws.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 3 now belongs to the new 6d2a7d87 file (was cbcdd771).
$ws2.Range("A3").Value = "6d2a7d87-f970-42ce-bf51-6d3c850207dd.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "False"
$ws2.Range("G3").Value = "6d2a7d87-f970-42ce-bf51-6d3c850207dd.197eded308f8ed22b21c9850f80743b3961c023e.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-21 00:47:49"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

# Row 4: cbcdd771 (moved down from row 3).
$ws2.Range("A4").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "e2e"
$ws2.Range("E4").Value = "ht"
$ws2.Range("F4").Value = "False"
$ws2.Range("G4").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.b8e4142af020d03b283755bd354fcda2d644bedb.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-21 00:46:47"
$ws2.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I4").Value = ""
$ws2.Range("J4").Value = ""
$ws2.Range("K4").Value = "0001-01-01 00:00:00"
$ws2.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L4").Value = ""
$ws2.Range("M4").Value = "True"
$ws2.Range("N4").Value = ""
$ws2.Range("O4").Value = "False"
$ws2.Range("P4").Value = ""

# Row 5: new edbab8bf file.
$ws2.Range("A5").Value = "edbab8bf-23e4-47e8-831b-14c9aee13e0e.md"
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("D5").Value = "e2e"
$ws2.Range("E5").Value = "ht"
$ws2.Range("F5").Value = "False"
$ws2.Range("G5").Value = "edbab8bf-23e4-47e8-831b-14c9aee13e0e.977c565182eb70f5592892745c773df547f5cf39.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-21 00:47:49"
$ws2.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I5").Value = ""
$ws2.Range("J5").Value = ""
$ws2.Range("K5").Value = "0001-01-01 00:00:00"
$ws2.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L5").Value = ""
$ws2.Range("M5").Value = "True"
$ws2.Range("N5").Value = ""
$ws2.Range("O5").Value = "False"
$ws2.Range("P5").Value = ""

# Rebuild hyperlinks (A2, I2, A3, A4, A5).
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a454c3bc198b02c2e49d4016dd6550af4a97165/e2e/98ad48ff-e5a0-4540-a055-8ec88ce579da.md", "", "", "98ad48ff-e5a0-4540-a055-8ec88ce579da.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/189ed7008381897ea17b8de5754d72ece0e5c0d0/e2e/98ad48ff-e5a0-4540-a055-8ec88ce579da.md", "", "", "98ad48ff-e5a0-4540-a055-8ec88ce579da.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197eded308f8ed22b21c9850f80743b3961c023e/e2e/6d2a7d87-f970-42ce-bf51-6d3c850207dd.md", "", "", "6d2a7d87-f970-42ce-bf51-6d3c850207dd.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb043c5240aa15b0c64d9bbc34b532dce3e332d3/e2e/cbcdd771-a78d-4d08-b66d-488e5b202a24.md", "", "", "cbcdd771-a78d-4d08-b66d-488e5b202a24.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/977c565182eb70f5592892745c773df547f5cf39/e2e/edbab8bf-23e4-47e8-831b-14c9aee13e0e.md", "", "", "edbab8bf-23e4-47e8-831b-14c9aee13e0e.md") | Out-Null

$loZhCn = $ws2.ListObjects.Item(1)
$loZhCn.Resize($ws2.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 3 now belongs to the new 6d2a7d87 file (was cbcdd771).
$ws3.Range("A3").Value = "6d2a7d87-f970-42ce-bf51-6d3c850207dd.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "False"
$ws3.Range("G3").Value = "6d2a7d87-f970-42ce-bf51-6d3c850207dd.197eded308f8ed22b21c9850f80743b3961c023e.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-21 00:47:53"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

# Row 4: cbcdd771 (moved down from row 3).
$ws3.Range("A4").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "e2e"
$ws3.Range("E4").Value = "ht"
$ws3.Range("F4").Value = "False"
$ws3.Range("G4").Value = "cbcdd771-a78d-4d08-b66d-488e5b202a24.b8e4142af020d03b283755bd354fcda2d644bedb.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-21 00:46:51"
$ws3.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I4").Value = ""
$ws3.Range("J4").Value = ""
$ws3.Range("K4").Value = "0001-01-01 00:00:00"
$ws3.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L4").Value = ""
$ws3.Range("M4").Value = "True"
$ws3.Range("N4").Value = ""
$ws3.Range("O4").Value = "False"
$ws3.Range("P4").Value = ""

# Row 5: new edbab8bf file.
$ws3.Range("A5").Value = "edbab8bf-23e4-47e8-831b-14c9aee13e0e.md"
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("D5").Value = "e2e"
$ws3.Range("E5").Value = "ht"
$ws3.Range("F5").Value = "False"
$ws3.Range("G5").Value = "edbab8bf-23e4-47e8-831b-14c9aee13e0e.977c565182eb70f5592892745c773df547f5cf39.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-21 00:47:53"
$ws3.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I5").Value = ""
$ws3.Range("J5").Value = ""
$ws3.Range("K5").Value = "0001-01-01 00:00:00"
$ws3.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L5").Value = ""
$ws3.Range("M5").Value = "True"
$ws3.Range("N5").Value = ""
$ws3.Range("O5").Value = "False"
$ws3.Range("P5").Value = ""

# Rebuild hyperlinks (A2, I2, A3, A4, A5).
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a454c3bc198b02c2e49d4016dd6550af4a97165/e2e/98ad48ff-e5a0-4540-a055-8ec88ce579da.md", "", "", "98ad48ff-e5a0-4540-a055-8ec88ce579da.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3c6a2e23770148e32a2fafe98b3df8f08acd87ac/e2e/98ad48ff-e5a0-4540-a055-8ec88ce579da.md", "", "", "98ad48ff-e5a0-4540-a055-8ec88ce579da.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/197eded308f8ed22b21c9850f80743b3961c023e/e2e/6d2a7d87-f970-42ce-bf51-6d3c850207dd.md", "", "", "6d2a7d87-f970-42ce-bf51-6d3c850207dd.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb043c5240aa15b0c64d9bbc34b532dce3e332d3/e2e/cbcdd771-a78d-4d08-b66d-488e5b202a24.md", "", "", "cbcdd771-a78d-4d08-b66d-488e5b202a24.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/977c565182eb70f5592892745c773df547f5cf39/e2e/edbab8bf-23e4-47e8-831b-14c9aee13e0e.md", "", "", "edbab8bf-23e4-47e8-831b-14c9aee13e0e.md") | Out-Null

$loDeDe = $ws3.ListObjects.Item(1)
$loDeDe.Resize($ws3.Range("A1:P5"))
